# Store_Template.xlsx: translate worksheet tab names from English to
# Vietnamese, nudge the saved cursor position on the "ParentStore" /
# "Đại lý cha" sheet from C9 to C10, and leave the workbook with the
# "Province" / "Tỉnh, Thành phố" sheet active (instead of the first
# sheet, "Store" / "Đại lý") the next time it is opened.

$wb = $excel.ActiveWorkbook

# --- Rename worksheet tabs (English -> Vietnamese) ---------------------
# Look sheets up by their current (pre-edit) English names so the order
# in which we rename them doesn't matter, then rename in place.
$wb.Worksheets.Item("Store").Name       = "Đại lý"
$wb.Worksheets.Item("Org").Name         = "Đơn vị tổ chức"
$wb.Worksheets.Item("ParentStore").Name = "Đại lý cha"
$wb.Worksheets.Item("StoreType").Name   = "Loại đại lý"
$wb.Worksheets.Item("StoreGroup").Name  = "Nhóm đại lý"
$wb.Worksheets.Item("Province").Name    = "Tỉnh, Thành phố"
$wb.Worksheets.Item("District").Name    = "Quận, Huyện"
$wb.Worksheets.Item("Ward").Name        = "Phường, Xã"
# "Quy tac import" keeps its original name.

# --- Update the saved selection on "Đại lý cha" (was C9, now C10) ------
$wsParentStore = $wb.Worksheets.Item("Đại lý cha")
$wsParentStore.Range("C10").Select() | Out-Null

# --- Make "Tỉnh, Thành phố" the active/selected tab on open -------------
# (previously the first sheet, "Đại lý", was the selected tab). Do this
# last so it "wins" as the tab that is active when the workbook re-opens.
$wsProvince = $wb.Worksheets.Item("Tỉnh, Thành phố")
$wsProvince.Activate()
